$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 2: replace the old DASH_SERV_01 (success) scenario with the new
# DASH_SERV_02 (exception / error-handling) scenario.
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "DASH_SERV_02"
$ws.Range("B2").Value = "Lỗi hệ thống (Exception)"
$ws.Range("C2").Value = "Error"
$ws.Range("D2").Value = "DAO ném lỗi RuntimeException"
$ws.Range("E2").Value = "Log Error & Forward JSP (Safe Mode)"
$ws.Range("F2").Value = "OK"
$ws.Range("G2").Value = "PASS"

# Re-style G2 from the old "FAIL" (red) look to the new "PASS" (green) look.
$ws.Range("G2").Font.Name = "Calibri"
$ws.Range("G2").Font.Bold = $true
$ws.Range("G2").Font.Color = 32768

# ---------------------------------------------------------------------------
# Row 3 (new row): the original DASH_SERV_01 scenario now passes, with a
# trimmed "Kết Quả Thực Tế" (actual result) of "OK" instead of the long
# Mockito failure dump, and a "PASS" status.
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "DASH_SERV_01"
$ws.Range("B3").Value = "Load trang Dashboard thành công"
$ws.Range("C3").Value = "Data: Rev=10tr, Orders=50"
$ws.Range("D3").Value = "1. Gọi DAO lấy số liệu`n2. Set attributes`n3. Forward JSP"
$ws.Range("E3").Value = "Forward Dashboard.jsp & Data OK"
$ws.Range("F3").Value = "OK"
$ws.Range("G3").Value = "PASS"

# Give G3 the same green "PASS" look by copying G2's now-correct format.
$ws.Range("G2").Copy()
$ws.Range("G3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Writing the multi-line "Các Bước" text auto-expanded the row heights;
# auto-fit puts them back to the (unexceptional) default row height so no
# stray row-height override is left behind, matching the original sheet.
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(3).AutoFit()

# ---------------------------------------------------------------------------
# Column widths: columns C-F were resized to fit the new (generally shorter)
# content now that the huge Mockito stack-trace text is gone from column F.
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 24.0
$ws.Columns.Item(4).ColumnWidth = 28.333333333333332
$ws.Columns.Item(5).ColumnWidth = 33.166666666666664
$ws.Columns.Item(6).ColumnWidth = 15.333333333333334

Write-Output "edit applied"
